$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = -8.826700000000002
$ws.Range("E3").Value = 16.06359999999999
$ws.Range("C12").Value = -10.7704
$ws.Range("D14").Value = -7.360100000000001
$ws.Range("E20").Value = 16.02309999999999
$ws.Range("E25").Value = 17.20420000000001
$ws.Range("D26").Value = -8.572199999999992
$ws.Range("C27").Value = -11.94860000000001
$ws.Range("E30").Value = 15.5898
$ws.Range("D31").Value = -9.293699999999991
$ws.Range("C32").Value = -12.4215
$ws.Range("D35").Value = -8.864899999999995
$ws.Range("C36").Value = -12.807
$ws.Range("D37").Value = -8.141399999999994
$ws.Range("C38").Value = -12.0965
$ws.Range("E44").Value = 16.12149999999999
$ws.Range("D45").Value = -7.384100000000002
$ws.Range("C46").Value = -14.50489999999999
$ws.Range("E47").Value = 16.121
$ws.Range("D52").Value = -7.670799999999997
$ws.Range("C54").Value = -13.08200000000001
$ws.Range("C55").Value = -13.42549999999999
$ws.Range("C56").Value = -12.35309999999999
$ws.Range("D57").Value = -8.498699999999998
$ws.Range("E58").Value = 16.5268
$ws.Range("C67").Value = -10.62460000000001
$ws.Range("C69").Value = -11.3826
$ws.Range("C72").Value = -11.42510000000001
$ws.Range("E78").Value = 16.54070000000002
$ws.Range("D81").Value = -7.102899999999998
$ws.Range("C83").Value = -14.03220000000001
$ws.Range("D83").Value = -8.783699999999996
$ws.Range("E84").Value = 16.57699999999999
$ws.Range("C86").Value = -13.57649999999999
$ws.Range("E89").Value = 17.36590000000002
$ws.Range("C91").Value = -10.40139999999999
$ws.Range("E91").Value = 17.93940000000002
$ws.Range("E92").Value = 18.02070000000002
$ws.Range("C93").Value = -11.1769
$ws.Range("E96").Value = 15.09569999999999
$ws.Range("C99").Value = -13.512
$ws.Range("D100").Value = -8.292200000000001
$ws.Range("D102").Value = -8.134400000000001
$ws.Range("E102").Value = 16.56210000000002
